$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text number format on the Price/Volume columns so that
# numeric-looking strings (e.g. "290.53") are stored as text, matching the
# original inline-string cell contents instead of being parsed as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '22.368.66'
$ws.Range("E2").Value = '  -4.49%  '

$ws.Range("D3").Value = '1.567.11'
$ws.Range("E3").Value = '  -4.72%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("E5").Value = '  -0.02%  '

$ws.Range("D6").Value = '290.53'
$ws.Range("E6").Value = '  -2.96%  '

$ws.Range("D7").Value = '0.3682'
$ws.Range("E7").Value = '  -2.86%  '

$ws.Range("D8").Value = '49.35'
$ws.Range("E8").Value = '  -1.11%  '

$ws.Range("D9").Value = '0.3393'
$ws.Range("E9").Value = '  -3.75%  '

$ws.Range("D10").Value = '1.169'
$ws.Range("E10").Value = '  -3.58%  '

$ws.Range("D11").Value = '0.07605'
$ws.Range("E11").Value = '  -5.79%  '

$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("D13").Value = '21.19'
$ws.Range("E13").Value = '  -3.85%  '

$ws.Range("E14").Value = '  -4.99%  '

$ws.Range("D15").Value = '6.894'
$ws.Range("E15").Value = '  -5.71%  '

$ws.Range("D16").Value = '1.579.40'
$ws.Range("E16").Value = '  -4.21%  '

$ws.Range("D17").Value = '0.00001137'
$ws.Range("E17").Value = '  -5.36%  '

$ws.Range("D18").Value = '89.11'
$ws.Range("E18").Value = '  -7.83%  '

$ws.Range("D19").Value = '0.06761'
$ws.Range("E19").Value = '  -3.31%  '

$ws.Range("D21").Value = '6.231'
$ws.Range("E21").Value = '  -7.32%  '

$ws.Range("D22").Value = '0.5336'
$ws.Range("E22").Value = '  -6.97%  '

$ws.Range("E23").Value = '  -4.87%  '

$ws.Range("D24").Value = '12.03'
$ws.Range("E24").Value = '  -2.66%  '

$ws.Range("D25").Value = '22.375.60'
$ws.Range("E25").Value = '  -4.51%  '

$ws.Range("D26").Value = '2.384'
$ws.Range("E26").Value = '  -4.59%  '

$ws.Range("D27").Value = '2.985'
$ws.Range("E27").Value = '  +2.71%  '

$ws.Range("E28").Value = '  -4.58%  '

$ws.Range("D29").Value = '145.61'
$ws.Range("E29").Value = '  -4.99%  '

$ws.Range("D30").Value = '4.951'
$ws.Range("E30").Value = '  -4.95%  '

$ws.Range("D31").Value = '125.35'
$ws.Range("E31").Value = '  -5.35%  '

$ws.Range("D32").Value = '1.751.56'
$ws.Range("E32").Value = '  -4.19%  '

$ws.Range("D33").Value = '1.037'
$ws.Range("E33").Value = '  +5.81%  '

$ws.Range("D34").Value = '6.252'
$ws.Range("E34").Value = '  -8.99%  '

$ws.Range("D35").Value = '1.990'
$ws.Range("E35").Value = '  -6.16%  '

$ws.Range("D36").Value = '10.29'
$ws.Range("E36").Value = '  -9.78%  '

$ws.Range("D37").Value = '0.08463'
$ws.Range("E37").Value = '  -3.27%  '

$ws.Range("D38").Value = '0.02540'
$ws.Range("E38").Value = '  -6.02%  '

$ws.Range("D39").Value = '0.2328'
$ws.Range("E39").Value = '  -4.22%  '

$ws.Range("D40").Value = '0.06561'
$ws.Range("E40").Value = '  -3.68%  '

$ws.Range("D41").Value = '5.527'
$ws.Range("E41").Value = '  -6.40%  '

$ws.Range("D42").Value = '11.84'
$ws.Range("E42").Value = '  -7.75%  '

$ws.Range("D43").Value = '1.245'
$ws.Range("E43").Value = '  -3.81%  '

$ws.Range("E44").Value = '  -7.22%  '

$ws.Range("D45").Value = '14.33'

$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").Value = '0.6009'
$ws.Range("E47").Value = '  -5.14%  '

$ws.Range("D48").Value = '3.781'

$ws.Range("D49").Value = '2.133'
$ws.Range("E49").Value = '  -5.35%  '

$ws.Range("D50").Value = '1.271'
$ws.Range("E50").Value = '  +8.06%  '

$ws.Range("D51").Value = '123.57'
$ws.Range("E51").Value = '  -2.73%  '

# Restore the original (default) cell style now that the text values are set,
# so the cells do not end up carrying an explicit Text style like the source file.
$dataRange.Style = "Normal"
